# ChainSight Module3Output - "sync rest transaction files"
# Adds the in-transit / inventory / fcst / order-AO reconciliation block
# below the existing NetDemand table and highlights the C816 rows plus
# the new helper cells in yellow, matching the author's manual edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Highlight the existing "80813644" (C816-related) rows 6-8 in yellow ---
$ws.Range("A6:H8").Interior.Color = 65535

# --- New reconciliation / check cells ---
# Row 10: difference between the two C816 order rows
$ws.Range("F10").Interior.Color = 65535
$ws.Range("F10").Formula = "=F8-F7"

# Row 13: standalone check value
$ws.Range("H13").Interior.Color = 65535
$ws.Range("H13").Value = 24

# Row 15: headers for the new mini inventory/in-transit/fcst/order block
# (written fcst/order-AO before inventory/in-transit so the shared-string
# table is appended in the same order the author's session produced it)
$ws.Range("F15").Value = "fcst"
$ws.Range("G15").Value = "order/AO"
$ws.Range("D15").Value = "inventory"
$ws.Range("E15").Value = "in-transit"

# Row 20: sample values + balance formula
$ws.Range("D20").Value = 50
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 497
$ws.Range("G20").Interior.Color = 65535
$ws.Range("G20").Formula = "=D20+E20-F20"

# --- Restore the active cell to where the author left off editing ---
$ws.Range("I5").Select()

Write-Output "applied sync rest transaction files edits"
